# Generate Report for Handback
# The e2e\ea69d92a-...md file has now been handed back (in sync with en-US),
# so update its status/handback-datetime across all three report sheets and
# clear the stale "handback not latest" error detail.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the ea69d92a-... file ---
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the ea69d92a-... file ---
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-31 16:56:37"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.74705287

# --- de-de sheet: row 3 is the ea69d92a-... file ---
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-31 16:56:44"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.74705287
